$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.102.80"
$ws.Range("E2").Value = "  +5.30%  "
$ws.Range("D3").Value = "1.880.13"
$ws.Range("E3").Value = "  +3.95%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.ClearFormats()
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "280.94"
$c.ClearFormats()
$ws.Range("E5").Value = "  +2.18%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9993"
$c.ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5294"
$c.ClearFormats()
$ws.Range("E7").Value = "  +4.91%  "
$ws.Range("E8").Value = "  +0.42%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "45.44"
$c.ClearFormats()
$ws.Range("E9").Value = "  +2.21%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07044"
$c.ClearFormats()
$ws.Range("E10").Value = "  +6.08%  "
$ws.Range("E11").Value = "  +1.92%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.8173"
$c.ClearFormats()
$ws.Range("E12").Value = "  -1.80%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.07812"
$c.ClearFormats()
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "1.885.24"
$ws.Range("E14").Value = "  +4.28%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.199"
$c.ClearFormats()
$ws.Range("E15").Value = "  +2.93%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "90.51"
$c.ClearFormats()
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  +5.06%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000008195"
$c.ClearFormats()
$ws.Range("E19").Value = "  +2.62%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.ClearFormats()
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "27.130.91"
$ws.Range("D22").Value = "2.120.32"
$ws.Range("E22").Value = "  +4.36%  "
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  +1.77%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "6.215"
$c.ClearFormats()
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("E26").Value = "  +12.10%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "146.35"
$c.ClearFormats()
$ws.Range("E27").Value = "  +3.36%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.57"
$c.ClearFormats()
$ws.Range("E28").Value = "  +3.73%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.672"
$c.ClearFormats()
$ws.Range("E29").Value = "  +1.27%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "112.51"
$c.ClearFormats()
$ws.Range("E30").Value = "  +3.69%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.397"
$c.ClearFormats()
$ws.Range("E31").Value = "  +1.65%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.381"
$c.ClearFormats()
$ws.Range("E32").Value = "  +4.54%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.08921"
$c.ClearFormats()
$ws.Range("E33").Value = "  +1.81%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04902"
$c.ClearFormats()
$ws.Range("E34").Value = "  +2.17%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.177"
$c.ClearFormats()
$ws.Range("E35").Value = "  +3.80%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.7469"
$c.ClearFormats()
$ws.Range("E36").Value = "  +3.21%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.899"
$c.ClearFormats()
$ws.Range("E37").Value = "  +0.71%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.306"
$c.ClearFormats()
$ws.Range("E38").Value = "  +9.12%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.419"
$c.ClearFormats()
$ws.Range("E39").Value = "  +6.49%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5325"
$c.ClearFormats()
$ws.Range("E40").Value = "  +3.05%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.01881"
$c.ClearFormats()
$ws.Range("E41").Value = "  +1.39%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.9723"
$c.ClearFormats()
$ws.Range("E42").Value = "  +2.81%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "117.07"
$c.ClearFormats()
$ws.Range("E43").Value = "  +4.30%  "
$ws.Range("E44").Value = "  +2.56%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "8.221"
$c.ClearFormats()
$ws.Range("E45").Value = "  +3.13%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.ClearFormats()
$ws.Range("E46").Value = "  -0.07%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.4597"
$c.ClearFormats()
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("E48").Value = "  -0.47%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.488"
$c.ClearFormats()
$ws.Range("E49").Value = "  +2.56%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "36.72"
$c.ClearFormats()
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.529"
$c.ClearFormats()
$ws.Range("E51").Value = "  +2.35%  "
